$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.155.50"
$ws.Range("E2").Value = "  +1.54%  "
$ws.Range("D3").Value = "3.892.82"
$ws.Range("E3").Value = "  -0.44%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").Value = "'525.88"
$ws.Range("E5").Value = "  +8.61%  "
$ws.Range("D6").Value = "'143.50"
$ws.Range("E6").Value = "  -1.68%  "
$ws.Range("D7").Value = "'0.608"
$ws.Range("E7").Value = "  -1.81%  "
$ws.Range("D8").Value = "'0.998"
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "'0.718"
$ws.Range("E9").Value = "  -1.77%  "
$ws.Range("D10").Value = "'0.170"
$ws.Range("E10").Value = "  +1.86%  "
$ws.Range("D11").Value = "'0.0000330"
$ws.Range("E11").Value = "  -4.77%  "
$ws.Range("D12").Value = "'41.91"
$ws.Range("E12").Value = "  -2.65%  "
$ws.Range("D13").Value = "4.498.03"
$ws.Range("E13").Value = "  -0.97%  "
$ws.Range("D14").Value = "'10.17"
$ws.Range("E14").Value = "  -4.36%  "
$ws.Range("D15").Value = "4.034.52"
$ws.Range("E15").Value = "  +3.52%  "
$ws.Range("E16").Value = "  -0.56%  "
$ws.Range("D17").Value = "'1.22"
$ws.Range("E17").Value = "  +7.64%  "
$ws.Range("D18").Value = "'13.74"
$ws.Range("E18").Value = "  -3.13%  "
$ws.Range("D19").Value = "'19.64"
$ws.Range("E19").Value = "  -2.02%  "
$ws.Range("D20").Value = "69.015.28"
$ws.Range("E20").Value = "  +1.13%  "
$ws.Range("D21").Value = "'423.70"
$ws.Range("E21").Value = "  -1.26%  "
$ws.Range("D22").Value = "'3.32"
$ws.Range("E22").Value = "  -4.82%  "
$ws.Range("D23").Value = "'14.10"
$ws.Range("E23").Value = "  -6.55%  "
$ws.Range("B24").Value = "PancakeSwap"
$ws.Range("C24").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D24").Value = "'4.07"
$ws.Range("E24").Value = "  +9.97%  "
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").Value = "'87.32"
$ws.Range("E25").Value = "  -0.97%  "
$ws.Range("D26").Value = "'11.51"
$ws.Range("E26").Value = "  -1.89%  "
$ws.Range("D27").Value = "'10.54"
$ws.Range("E27").Value = "  -5.53%  "
$ws.Range("D28").Value = "'35.83"
$ws.Range("E28").Value = "  -4.75%  "
$ws.Range("D29").Value = "'695.43"
$ws.Range("E29").Value = "  -3.00%  "
$ws.Range("D30").Value = "'13.08"
$ws.Range("E30").Value = "  -5.13%  "
$ws.Range("D31").Value = "'0.125"
$ws.Range("E31").Value = "  -4.36%  "
$ws.Range("D32").Value = "'2.80"
$ws.Range("E32").Value = "  -3.94%  "
$ws.Range("D33").Value = "'67.79"
$ws.Range("E33").Value = "  +11.10%  "
$ws.Range("D34").Value = "'0.443"
$ws.Range("E34").Value = "  +11.75%  "
$ws.Range("D35").Value = "'5.95"
$ws.Range("E35").Value = "  -3.37%  "
$ws.Range("D36").Value = "'40.05"
$ws.Range("E36").Value = "  -3.45%  "
$ws.Range("D37").Value = "0.0₃0842"
$ws.Range("E37").Value = "  -7.31%  "
$ws.Range("D38").Value = "'0.997"
$ws.Range("E38").Value = "  -0.14%  "
$ws.Range("D39").Value = "'0.147"
$ws.Range("E39").Value = "  +0.92%  "
$ws.Range("D40").Value = "'0.997"
$ws.Range("E40").Value = "  -0.50%  "
$ws.Range("D41").Value = "'0.0477"
$ws.Range("E41").Value = "  -2.61%  "
$ws.Range("B42").Value = "Fetch.AI"
$ws.Range("C42").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D42").Value = "'2.78"
$ws.Range("E42").Value = "  -7.75%  "
$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").Value = "'3.00"
$ws.Range("E43").Value = "  -0.58%  "
$ws.Range("D44").Value = "'2.93"
$ws.Range("E44").Value = "  -6.53%  "
$ws.Range("E45").Value = "  -0.71%  "
$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").Value = "'0.139"
$ws.Range("E46").Value = "  -1.62%  "
$ws.Range("B47").Value = "Stacks"
$ws.Range("C47").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D47").Value = "'3.04"
$ws.Range("E47").Value = "  +8.65%  "
$ws.Range("B48").Value = "LidoDAOToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D48").Value = "'3.27"
$ws.Range("E48").Value = "  -4.25%  "
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").Value = "'142.08"
$ws.Range("E49").Value = "  -1.91%  "
$ws.Range("B50").Value = "ARBITRUM"
$ws.Range("C50").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D50").Value = "'2.04"
$ws.Range("E50").Value = "  -4.44%  "
$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").Value = "0.0₆0326"
$ws.Range("E51").Value = "  -3.42%  "
